$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.992.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.356.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.32"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.64%  "
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.107"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.708.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.901"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.354.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.929.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.52%  "
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "259.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +24.15%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "176.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("E36").Value = "  +5.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.08%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("E41").Value = "  +12.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.205"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.38%  "
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.26%  "
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.78%  "
